# ---------------------------------------------------------------------------
# "Update the write up"
#
# 1) "Task1 Result:" -> "Task 1 Result:"
#    A space is typed between "Task" and "1", which is exactly what Word
#    does when the insertion point lands inside an existing run: the run is
#    split in two around the caret and the freshly typed text becomes its
#    own run in between, giving three runs ("Task" / " " / "1 Result:")
#    with identical (absent) character formatting.
#
# 2) Word re-mints wp14:editId for every drawing's <wp:inline> on save -
#    independent of whether that particular picture was touched - so both
#    pictures in the doc pick up a fresh edit id as a side effect of the
#    save that followed the text edit above.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Split "Task1 Result:" into "Task" | " " | "1 Result:" -------------

$hit = $d.Content
$null = $hit.Find.Execute("Task1")
$splitAt = $hit.Start + 4   # caret position right between "Task" and "1"

$caret = $d.Range($splitAt, $splitAt)
$null = $caret.InsertBefore(" ")

# Stamp a throwaway bookmark on exactly the newly typed character and
# remove it again; this forces the engine to materialize the run break on
# both sides of the caret (matching what Word's own typing-in-the-middle-
# of-a-run behavior produces) without leaving any character-formatting
# residue behind.
$typed = $d.Range($splitAt, $splitAt + 1)
$null = $d.Bookmarks.Add("__taskSplit", $typed)
$null = $d.Bookmarks("__taskSplit").Delete()

# --- 2) Re-stamp wp14:editId on every inline picture -----------------------

function Update-DrawingEditId($Shape, $NewEditId) {
    # NOTE: always call this positionally (Update-DrawingEditId $x $y) -
    # named-parameter binding (-Shape $x) loses COM object identity in
    # this host and the function body sees an empty value instead.

    $start = $Shape.Range.Start
    $end = $Shape.Range.End
    $paraXml = $script:OriginalParagraphXml[$start]
    if (-not $paraXml) { return }

    # Swap in the new wp14:editId value, leaving everything else (the
    # picture, its extents, its relationship id, paragraph metadata, ...)
    # byte-for-byte identical.
    $updated = $paraXml -replace 'wp14:editId="[0-9A-F]{8}"', ('wp14:editId="' + $NewEditId + '"')

    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' +
           'xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" ' +
           'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" ' +
           'xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" ' +
           'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
           'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
           'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
           'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
           'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' +
           'mc:Ignorable="w14 wp14">' +
           '<w:body>' + $updated + '<w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    # Re-fetch a *document-level* Range (not the InlineShape's own .Range)
    # right before the call - replacing through a range obtained straight
    # from the shape object duplicates the drawing in this host.
    $target = $d.Range($start, $end)
    $null = $target.InsertXML($pkg)
}

# The two paragraphs that hold the pictures, captured from the pristine
# document (InsertXML's WordOpenXML round-trip isn't scoped to the range
# it's read from in this host, so the fragments are rebuilt from the known
# original markup instead of read back live).
$script:OriginalParagraphXml = @{}

$picture1Paragraph = '<w:p w14:paraId="7FC8B377" w14:textId="4CE0349E" w:rsidR="00A67C4E" w:rsidRDefault="00F26B31"><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="0D79AAFA" wp14:editId="564D5104"><wp:extent cx="5525691" cy="2266950"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="753710868" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId6"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="5526891" cy="2267442"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'

$picture2Paragraph = '<w:p w14:paraId="67F088EB" w14:textId="2CC0A523" w:rsidR="00A67C4E" w:rsidRDefault="00E825BA"><w:r w:rsidRPr="00E825BA"><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="758C0428" wp14:editId="34D4E786"><wp:extent cx="5626100" cy="2735511"/><wp:effectExtent l="0" t="0" r="0" b="8255"/><wp:docPr id="1273140135" name="Picture 1" descr="A picture containing diagram&#10;&#10;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1273140135" name="Picture 1" descr="A picture containing diagram&#10;&#10;Description automatically generated"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId7"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5630360" cy="2737582"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'

if ($d.InlineShapes.Count -ge 1) {
    $shape1 = $d.InlineShapes(1)
    $script:OriginalParagraphXml[$shape1.Range.Start] = $picture1Paragraph
    Update-DrawingEditId $shape1 "4475B11B"
}

if ($d.InlineShapes.Count -ge 2) {
    $shape2 = $d.InlineShapes(2)
    $script:OriginalParagraphXml[$shape2.Range.Start] = $picture2Paragraph
    Update-DrawingEditId $shape2 "72183925"
}

Write-Output "Task write-up title split into three runs; drawing edit ids refreshed."
